$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Cells whose new value looks like a pure number need to be forced to remain
# text (matching the original inlineStr/shared-string text representation),
# otherwise Excel auto-converts them to floating point numbers.
$textForceCells = @(
    @{ Cell = "D5"; Value = "594.01" }
    @{ Cell = "D6"; Value = "179.06" }
    @{ Cell = "D14"; Value = "26.89" }
    @{ Cell = "D21"; Value = "367.70" }
    @{ Cell = "D26"; Value = "70.80" }
    @{ Cell = "D29"; Value = "0.996" }
    @{ Cell = "D31"; Value = "541.71" }
    @{ Cell = "D32"; Value = "8.31" }
    @{ Cell = "D38"; Value = "157.42" }
    @{ Cell = "D39"; Value = "18.88" }
    @{ Cell = "D41"; Value = "0.357" }
    @{ Cell = "D43"; Value = "5.23" }
    @{ Cell = "D44"; Value = "2.58" }
    @{ Cell = "D46"; Value = "147.64" }
    @{ Cell = "D47"; Value = "0.563" }
)
foreach ($item in $textForceCells) {
    $rng = $ws.Range($item.Cell)
    $rng.NumberFormat = "@"
    $rng.Value = $item.Value
    $rng.Style = "Normal"
}

# Remaining cells: values are safely non-numeric text already.
$ws.Range("D2").Value = "68.234.23"
$ws.Range("E2").Value = "  +2.56%  "
$ws.Range("D3").Value = "2.533.70"
$ws.Range("E3").Value = "  +0.83%  "
$ws.Range("E4").Value = "  -0.01%  "
$ws.Range("E5").Value = "  +1.86%  "
$ws.Range("E6").Value = "  +4.62%  "
$ws.Range("E7").Value = "  -0.06%  "
$ws.Range("E8").Value = "  +1.16%  "
$ws.Range("D9").Value = "2.531.95"
$ws.Range("E9").Value = "  +0.81%  "
$ws.Range("E10").Value = "  +2.58%  "
$ws.Range("E12").Value = "  -0.07%  "
$ws.Range("E13").Value = "  -0.40%  "
$ws.Range("E14").Value = "  +0.63%  "
$ws.Range("D15").Value = "2.995.00"
$ws.Range("E15").Value = "  +0.84%  "
$ws.Range("E16").Value = "  +2.30%  "
$ws.Range("D17").Value = "68.141.84"
$ws.Range("E17").Value = "  +2.76%  "
$ws.Range("D18").Value = "2.531.87"
$ws.Range("E18").Value = "  +0.53%  "
$ws.Range("E19").Value = "  +2.42%  "
$ws.Range("E20").Value = "  +2.36%  "
$ws.Range("E21").Value = "  +5.66%  "
$ws.Range("E22").Value = "  +0.15%  "
$ws.Range("E23").Value = "  +1.70%  "
$ws.Range("E24").Value = "  -1.48%  "
$ws.Range("E25").Value = "  -0.04%  "
$ws.Range("E26").Value = "  +1.17%  "
$ws.Range("E27").Value = "  +3.03%  "
$ws.Range("E29").Value = "  -0.44%  "
$ws.Range("D30").Value = "0.0₂01000"
$ws.Range("E30").Value = "  +2.51%  "
$ws.Range("E31").Value = "  +3.36%  "
$ws.Range("E32").Value = "  +2.69%  "
$ws.Range("E33").Value = "  +1.95%  "
$ws.Range("E35").Value = "  -0.78%  "
$ws.Range("E36").Value = "  +0.07%  "
$ws.Range("E38").Value = "  +0.13%  "
$ws.Range("E40").Value = "  +1.68%  "
$ws.Range("E41").Value = "  +0.40%  "
$ws.Range("E43").Value = "  +2.72%  "
$ws.Range("E44").Value = "  +3.68%  "
$ws.Range("E45").Value = "  -0.05%  "
$ws.Range("E46").Value = "  -0.94%  "
$ws.Range("E47").Value = "  +0.81%  "
$ws.Range("D49").Value = "0.0₆0279"
$ws.Range("E49").Value = "  +3.51%  "
$ws.Range("E50").Value = "  -0.44%  "
$ws.Range("E51").Value = "  +0.36%  "
